$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for a new item row by inserting a blank row just above the
#    "totals" row (row 17), which pushes the totals row (17 -> 18) and the
#    footer row (18 -> 19) down by one without touching the existing item
#    rows 7-16.
# ---------------------------------------------------------------------------
$ws.Rows("17:17").Insert()

# Copy formatting (styles) only - from the last item row (16) into the new
# row (17) so it looks just like the other item rows.
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights: new item row matches the other item rows; totals row keeps
# its original (slightly smaller) height; footer keeps its height.
$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("18:18").RowHeight = 24.75

# Recreate the merged cells for the new item row.
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# ---------------------------------------------------------------------------
# 2) Shift the existing item rows 9-16 down into 10-17 (bottom-up so that
#    values are not overwritten before being read), freeing up row 9 for
#    the newly sold/reported item.
# ---------------------------------------------------------------------------
for ($r = 16; $r -ge 9; $r--) {
    $rSrc = $r
    $rDst = $r + 1

    $ws.Range("C$rDst").Value2 = $ws.Range("C$rSrc").Value2
    $ws.Range("H$rDst").Value2 = $ws.Range("H$rSrc").Value2
    $ws.Range("N$rDst").Value2 = $ws.Range("N$rSrc").Value2
    $ws.Range("Q$rDst").Value2 = $ws.Range("Q$rSrc").Value2

    # L and P columns use a numeric display format, but the source file
    # stores their contents as plain text - force text formatting while
    # assigning the value, then restore the original numeric format so the
    # cell keeps its original style / appearance.
    $ws.Range("L$rDst").NumberFormat = "@"
    $ws.Range("L$rDst").Value2 = $ws.Range("L$rSrc").Value2
    $ws.Range("L$rDst").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

    $ws.Range("P$rDst").NumberFormat = "@"
    $ws.Range("P$rDst").Value2 = $ws.Range("P$rSrc").Value2
    $ws.Range("P$rDst").NumberFormat = "0.00"
}

# ---------------------------------------------------------------------------
# 3) Fill the freed-up row 9 with the new item that was reported in this
#    version of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("C9").Value2 = "DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP."
$ws.Range("H9").Value2 = "0:1"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value2 = "1"
$ws.Range("L9").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N9").Value2 = "36.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value2 = "11.8800"
$ws.Range("P9").NumberFormat = "0.00"

$ws.Range("Q9").Value2 = "0:1"

# The new last item row is number 11 in the sequential numbering column.
$ws.Range("A17").Value2 = 11

# ---------------------------------------------------------------------------
# 4) Update the totals row (now row 18) and the footer timestamp (now in
#    row 19) to reflect the added item.
# ---------------------------------------------------------------------------
$ws.Range("P18").Value2 = 357.19999999999999
$ws.Range("A19").Value2 = "Thursday, 7 August, 2025 10:29 AM"
